$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33: Glazed and Confused
$ws.Range("H33").Value = 742
$ws.Range("I33").Value = 722.93335
$ws.Range("J33").Value = 782.8570999999999
$ws.Range("K33").Value = 722.93335
$ws.Range("L33").Value = 782.8570999999999
$ws.Range("M33").Value = -493.93335
$ws.Range("N33").Value = -1240.8571
# Row 62: The Mustache Suits Him
$ws.Range("H62").Value = 19177.928
$ws.Range("I62").Value = 36781.832
$ws.Range("K62").Value = 36781.832
$ws.Range("M62").Value = -36157.832
# Row 65: Forgery of Convenience (L)
$ws.Range("H65").Value = 19177.928
$ws.Range("I65").Value = 36781.832
$ws.Range("K65").Value = 183909.16
$ws.Range("M65").Value = -180789.16
# Row 116: Growing Up
$ws.Range("H116").Value = 38310.535
$ws.Range("I116").Value = 79023.28999999999
$ws.Range("J116").Value = 2686.875
$ws.Range("K116").Value = 79023.28999999999
$ws.Range("L116").Value = 2686.875
$ws.Range("M116").Value = -75581.28999999999
$ws.Range("N116").Value = -9570.875
# Row 125: Body over Mind
$ws.Range("H125").Value = 83333990
$ws.Range("I125").Value = 475
$ws.Range("K125").Value = 4275
$ws.Range("M125").Value = -1815
# Row 132: Fast-forwarding Flora
$ws.Range("H132").Value = 4294.25
$ws.Range("I132").Value = 2778.8
$ws.Range("J132").Value = 6820
$ws.Range("K132").Value = 8336.400000000001
$ws.Range("L132").Value = 20460
$ws.Range("M132").Value = -5806.400000000001
$ws.Range("N132").Value = -25520

$ws = $wb.Worksheets.Item("ARM")
# Row 2: Ain't Got No Ingots
$ws.Range("H2").Value = 5124
$ws.Range("I2").Value = 3264.5
$ws.Range("J2").Value = 20000
$ws.Range("K2").Value = 3264.5
$ws.Range("L2").Value = 20000
$ws.Range("M2").Value = -3151.5
$ws.Range("N2").Value = -20226
# Row 32: Ingot We Trust
$ws.Range("H32").Value = 2803.24
$ws.Range("I32").Value = 1842.1609
$ws.Range("J32").Value = 9235.076999999999
$ws.Range("K32").Value = 1842.1609
$ws.Range("L32").Value = 9235.076999999999
$ws.Range("M32").Value = -1555.1609
$ws.Range("N32").Value = -9809.076999999999
# Row 45: Hollow Hallmarks
$ws.Range("H45").Value = 1262.1818
$ws.Range("I45").Value = 1274
$ws.Range("J45").Value = 1252.3334
$ws.Range("K45").Value = 1274
$ws.Range("L45").Value = 1252.3334
$ws.Range("M45").Value = -897
$ws.Range("N45").Value = -2006.3334
# Row 61: Dealing with the Tough Stuff
$ws.Range("H61").Value = 422641.2
$ws.Range("I61").Value = 368593.47
$ws.Range("J61").Value = 503712.8
$ws.Range("K61").Value = 368593.47
$ws.Range("L61").Value = 503712.8
$ws.Range("M61").Value = -368381.47
$ws.Range("N61").Value = -504136.8
# Row 74: As the Bolt Flies
$ws.Range("H74").Value = 201677.4
$ws.Range("I74").Value = 228278.27
$ws.Range("J74").Value = 95273.91
$ws.Range("K74").Value = 228278.27
$ws.Range("L74").Value = 95273.91
$ws.Range("M74").Value = -227404.27
$ws.Range("N74").Value = -97021.91
# Row 77: Heavy Metal Banned (L)
$ws.Range("H77").Value = 201677.4
$ws.Range("I77").Value = 228278.27
$ws.Range("J77").Value = 95273.91
$ws.Range("K77").Value = 1141391.35
$ws.Range("L77").Value = 476369.55
$ws.Range("M77").Value = -1137023.35
$ws.Range("N77").Value = -485105.55
# Row 96: The Gauntlet Is Cast
$ws.Range("H96").Value = 19896
$ws.Range("J96").Value = 19896
$ws.Range("L96").Value = 19896
$ws.Range("N96").Value = -25388
# Row 110: Scheduled Maintenance
$ws.Range("H110").Value = 1412.2106
$ws.Range("I110").Value = 1623.7142
$ws.Range("J110").Value = 820
$ws.Range("K110").Value = 1623.7142
$ws.Range("L110").Value = 820
$ws.Range("M110").Value = 421.2858000000001
$ws.Range("N110").Value = -4910
# Row 116: No Scope
$ws.Range("H116").Value = 5124
$ws.Range("I116").Value = 3264.5
$ws.Range("J116").Value = 20000
$ws.Range("K116").Value = 3264.5
$ws.Range("L116").Value = 20000
$ws.Range("M116").Value = -970.5
$ws.Range("N116").Value = -24588
# Row 122: Haste for High Durium
$ws.Range("H122").Value = 4100.725
$ws.Range("I122").Value = 4175.839
$ws.Range("K122").Value = 12527.517
$ws.Range("M122").Value = -10077.517
# Row 136: Metal with Mettle
$ws.Range("H136").Value = 422641.2
$ws.Range("I136").Value = 368593.47
$ws.Range("J136").Value = 503712.8
$ws.Range("K136").Value = 1105780.41
$ws.Range("L136").Value = 1511138.4
$ws.Range("M136").Value = -1103230.41
$ws.Range("N136").Value = -1516238.4

$ws = $wb.Worksheets.Item("BSM")
# Row 3: Hells Bells
$ws.Range("H3").Value = 5124
$ws.Range("I3").Value = 3264.5
$ws.Range("J3").Value = 20000
$ws.Range("K3").Value = 3264.5
$ws.Range("L3").Value = 20000
$ws.Range("M3").Value = -3150.5
$ws.Range("N3").Value = -20228
# Row 132: Always Be Prepaired
$ws.Range("H132").Value = 27890
$ws.Range("J132").Value = 27890
$ws.Range("L132").Value = 27890
$ws.Range("N132").Value = -38010
# Row 134: Ruthenium Supremium
$ws.Range("H134").Value = 3936.1738
$ws.Range("I134").Value = 3997.7273
$ws.Range("J134").Value = 3879.75
$ws.Range("K134").Value = 11993.1819
$ws.Range("L134").Value = 11639.25
$ws.Range("M134").Value = -9458.1819
$ws.Range("N134").Value = -16709.25

$ws = $wb.Worksheets.Item("CRP")
# Row 122: Timber of Tenkonto
$ws.Range("H122").Value = 1900.9166
$ws.Range("I122").Value = 1070.6666
$ws.Range("K122").Value = 3211.9998
$ws.Range("M122").Value = -761.9998000000001
# Row 134: Wood You Be Quiet
$ws.Range("H134").Value = 2015.4166
$ws.Range("I134").Value = 1293.9412
$ws.Range("J134").Value = 2660.9473
$ws.Range("K134").Value = 3881.8236
$ws.Range("L134").Value = 7982.841899999999
$ws.Range("M134").Value = -1346.8236
$ws.Range("N134").Value = -13052.8419

$ws = $wb.Worksheets.Item("CUL")
# Row 106: Herky Jerky
$ws.Range("H106").Value = 4275
$ws.Range("J106").Value = 4275
$ws.Range("L106").Value = 12825
$ws.Range("N106").Value = -14717

$ws = $wb.Worksheets.Item("GSM")
# Row 102: Put the Metal to the Peddle
$ws.Range("H102").Value = 5031.3335
$ws.Range("I102").Value = 2498.8572
$ws.Range("K102").Value = 2498.8572
$ws.Range("M102").Value = -876.8571999999999
# Row 122: Awarding Academic Excellence
$ws.Range("H122").Value = 1170.1
$ws.Range("I122").Value = 1225.125
$ws.Range("J122").Value = 950
$ws.Range("K122").Value = 3675.375
$ws.Range("L122").Value = 2850
$ws.Range("M122").Value = -1225.375
$ws.Range("N122").Value = -7750
# Row 125: Pewter-hewn Punishment
$ws.Range("H125").Value = 21333.334
$ws.Range("J125").Value = 21333.334
$ws.Range("L125").Value = 21333.334
$ws.Range("N125").Value = -26253.334

$ws = $wb.Worksheets.Item("LTW")
# Row 7: Tan Before the Ban
$ws.Range("H7").Value = 3291
$ws.Range("I7").Value = 3023.75
$ws.Range("J7").Value = 3504.8
$ws.Range("K7").Value = 3023.75
$ws.Range("L7").Value = 3504.8
$ws.Range("M7").Value = -2911.75
$ws.Range("N7").Value = -3728.8
# Row 40: Best Served Toad
$ws.Range("H40").Value = 3230.4546
$ws.Range("I40").Value = 3058.889
$ws.Range("J40").Value = 4002.5
$ws.Range("K40").Value = 3058.889
$ws.Range("L40").Value = 4002.5
$ws.Range("M40").Value = -2922.889
$ws.Range("N40").Value = -4274.5
# Row 55: It's Not a Job, It's a Calling
$ws.Range("H55").Value = 277.95456
$ws.Range("I55").Value = 167.61539
$ws.Range("J55").Value = 437.33334
$ws.Range("K55").Value = 167.61539
$ws.Range("L55").Value = 437.33334
$ws.Range("M55").Value = 5.384610000000009
$ws.Range("N55").Value = -783.33334
# Row 122: Hell on Leather
$ws.Range("H122").Value = 2547.4546
$ws.Range("I122").Value = 2253.6667
$ws.Range("J122").Value = 2900
$ws.Range("K122").Value = 6761.000100000001
$ws.Range("L122").Value = 8700
$ws.Range("M122").Value = -4311.000100000001
$ws.Range("N122").Value = -13600
# Row 126: Battered Books
$ws.Range("H126").Value = 3291
$ws.Range("I126").Value = 3023.75
$ws.Range("J126").Value = 3504.8
$ws.Range("K126").Value = 9071.25
$ws.Range("L126").Value = 10514.4
$ws.Range("M126").Value = -6601.25
$ws.Range("N126").Value = -15454.4
# Row 136: Respect for Br'aax
$ws.Range("H136").Value = 2923.6338
$ws.Range("I136").Value = 1457.0962
$ws.Range("J136").Value = 6937.316
$ws.Range("K136").Value = 4371.2886
$ws.Range("L136").Value = 20811.948
$ws.Range("M136").Value = -1821.2886
$ws.Range("N136").Value = -25911.948

$ws = $wb.Worksheets.Item("WVR")
# Row 122: Heavy Armoire
$ws.Range("H122").Value = 145144.86
$ws.Range("I122").Value = 334769.66
$ws.Range("J122").Value = 2926.25
$ws.Range("K122").Value = 1004308.98
$ws.Range("L122").Value = 8778.75
$ws.Range("M122").Value = -1001858.98
$ws.Range("N122").Value = -13678.75
# Row 126: A Polished Purchase
$ws.Range("H126").Value = 2327.7666
$ws.Range("I126").Value = 2542.25
$ws.Range("J126").Value = 2082.6428
$ws.Range("K126").Value = 7626.75
$ws.Range("L126").Value = 6247.928400000001
$ws.Range("M126").Value = -5156.75
$ws.Range("N126").Value = -11187.9284
